$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("E5").Value = 69
$ws.Range("F5").Value = 38
$ws.Range("H5").Value = 38

# Row 6
$ws.Range("E6").Value = 23

# Row 8
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 2
$ws.Range("H8").Value = 2

# Row 11
$ws.Range("E11").Value = 151

# Row 12
$ws.Range("E12").Value = 224
$ws.Range("F12").Value = 113
$ws.Range("H12").Value = 113

# Row 13
$ws.Range("E13").Value = 77

# Row 14
$ws.Range("E14").Value = 68

# Row 24
$ws.Range("E24").Value = 102
$ws.Range("F24").Value = 49
$ws.Range("H24").Value = 49

# Row 26
$ws.Range("E26").Value = 60

# Row 27
$ws.Range("E27").Value = 145

# Row 29
$ws.Range("E29").Value = 92

# Row 30
$ws.Range("E30").Value = 103

# Row 33
$ws.Range("E33").Value = 132
$ws.Range("F33").Value = 61
$ws.Range("H33").Value = 61

# Row 34
$ws.Range("E34").Value = 103

# Row 39
$ws.Range("E39").Value = 113

# Row 40
$ws.Range("E40").Value = 140

# Row 41
$ws.Range("E41").Value = 178
$ws.Range("F41").Value = 62
$ws.Range("H41").Value = 62

# Row 42
$ws.Range("E42").Value = 154
$ws.Range("F42").Value = 72
$ws.Range("H42").Value = 72

# Row 43
$ws.Range("E43").Value = 48

# Row 44
$ws.Range("E44").Value = 142

# Row 45
$ws.Range("E45").Value = 59

# Row 46
$ws.Range("E46").Value = 123

# Row 47
$ws.Range("E47").Value = 206

# Row 48
$ws.Range("E48").Value = 103
